$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 2")
$ws.Range("A2").Value = "chau"
